# Update "想去人数" (number of people interested) values on both the
# "展览" and "全部类型" sheets, which contain the same data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of cell -> new value
$updates = @{
    "F2"  = 1568
    "F7"  = 2630
    "F9"  = 1658
    "F12" = 544
    "F14" = 7
    "F15" = 64
    "F17" = 8
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cell in $updates.Keys) {
        $ws.Range($cell).Value = $updates[$cell]
    }
}
